$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1014688.4
$ws.Range("I33").Value = 1149673.8
$ws.Range("K33").Value = 1149673.8
$ws.Range("M33").Value = -1149444.8

$ws.Range("H40").Value = 2123
$ws.Range("J40").Value = 2612.5
$ws.Range("L40").Value = 2612.5
$ws.Range("N40").Value = -2962.5

$ws.Range("I104").Value = 608.6
$ws.Range("J104").Value = 1999
$ws.Range("K104").Value = 1825.8
$ws.Range("L104").Value = 5997
$ws.Range("M104").Value = -78.80000000000018
$ws.Range("N104").Value = -9491

$ws.Range("H137").Value = 57251.055
$ws.Range("I137").Value = 1817.3846
$ws.Range("K137").Value = 5452.1538
$ws.Range("M137").Value = -2902.1538

$ws.Range("H138").Value = 2396.923
$ws.Range("J138").Value = 5163.375
$ws.Range("L138").Value = 15490.125
$ws.Range("N138").Value = -25770.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 181245.33
$ws.Range("I32").Value = 177946.17
$ws.Range("K32").Value = 177946.17
$ws.Range("M32").Value = -177659.17

$ws.Range("H74").Value = 1849.4286
$ws.Range("J74").Value = 2568.8333
$ws.Range("L74").Value = 2568.8333
$ws.Range("N74").Value = -4316.8333

$ws.Range("H77").Value = 1849.4286
$ws.Range("J77").Value = 2568.8333
$ws.Range("L77").Value = 12844.1665
$ws.Range("N77").Value = -21580.1665

$ws.Range("H132").Value = 1640.4117
$ws.Range("I132").Value = 1492.9375
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 4478.8125
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -1948.8125
$ws.Range("N132").Value = -17060

$ws.Range("H139").Value = 71481.336
$ws.Range("J139").Value = 71481.336
$ws.Range("L139").Value = 71481.336
$ws.Range("N139").Value = -81761.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("N35").Value = 0

$ws.Range("H86").Value = 2219.5557
$ws.Range("I86").Value = 2152.7144
$ws.Range("J86").Value = 2453.5
$ws.Range("K86").Value = 2152.7144
$ws.Range("L86").Value = 2453.5
$ws.Range("M86").Value = -1029.7144
$ws.Range("N86").Value = -4699.5

$ws.Range("H89").Value = 2219.5557
$ws.Range("I89").Value = 2152.7144
$ws.Range("J89").Value = 2453.5
$ws.Range("K89").Value = 10763.572
$ws.Range("L89").Value = 12267.5
$ws.Range("M89").Value = -5147.572
$ws.Range("N89").Value = -23499.5

$ws.Range("H105").Value = 6668547
$ws.Range("I105").Value = 12502281
$ws.Range("K105").Value = 12502281
$ws.Range("M105").Value = -12500534

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").ClearContents()
$ws.Range("N51").Value = 0

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("N61").Value = 0

$ws.Range("H86").Value = 6927.6113
$ws.Range("I86").Value = 7035.3
$ws.Range("K86").Value = 7035.3
$ws.Range("M86").Value = -5912.3

$ws.Range("H89").Value = 6927.6113
$ws.Range("I89").Value = 7035.3
$ws.Range("K89").Value = 35176.5
$ws.Range("M89").Value = -29560.5

$ws.Range("H105").Value = 3254.5
$ws.Range("I105").Value = 3254.5
$ws.Range("K105").Value = 3254.5
$ws.Range("M105").Value = -1507.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 742.9
$ws.Range("I8").Value = 742.9
$ws.Range("K8").Value = 2228.7
$ws.Range("M8").Value = -2089.7

$ws.Range("H40").Value = 3821.875
$ws.Range("I40").Value = 79.166664
$ws.Range("J40").Value = 15050
$ws.Range("K40").Value = 316.666656
$ws.Range("L40").Value = 60200
$ws.Range("M40").Value = -247.666656
$ws.Range("N40").Value = -60338

$ws.Range("H63").Value = 300
$ws.Range("I63").Value = 300
$ws.Range("K63").Value = 900
$ws.Range("M63").Value = -151

$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 1000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2730

$ws.Range("H66").Value = 300
$ws.Range("I66").Value = 300
$ws.Range("K66").Value = 2700
$ws.Range("M66").Value = 1044

$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 1000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2064

$ws.Range("H70").Value = 2539
$ws.Range("I70").Value = 898.3333
$ws.Range("K70").Value = 2694.9999
$ws.Range("M70").Value = -2379.9999

$ws.Range("H73").Value = 2539
$ws.Range("I73").Value = 898.3333
$ws.Range("K73").Value = 2694.9999
$ws.Range("M73").Value = -1602.9999

$ws.Range("H121").Value = 27780076
$ws.Range("J121").Value = 3070.625
$ws.Range("L121").Value = 9211.875
$ws.Range("N121").Value = -11831.875

$ws.Range("H137").Value = 3800.8696
$ws.Range("I137").Value = 1791
$ws.Range("J137").Value = 4872.8
$ws.Range("K137").Value = 5373
$ws.Range("L137").Value = 14618.4
$ws.Range("M137").Value = -273
$ws.Range("N137").Value = -24818.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 20086.5
$ws.Range("I34").Value = 15000
$ws.Range("K34").Value = 15000
$ws.Range("M34").Value = -14732

$ws.Range("H76").Value = 20086.5
$ws.Range("I76").Value = 15000
$ws.Range("K76").Value = 15000
$ws.Range("M76").Value = -14685

$ws.Range("H79").Value = 20086.5
$ws.Range("I79").Value = 15000
$ws.Range("K79").Value = 15000
$ws.Range("M79").Value = -13908

$ws.Range("H125").Value = 75000
$ws.Range("J125").Value = 75000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -79920

$ws.Range("H132").Value = 2051.4614
$ws.Range("I132").Value = 1974.5454
$ws.Range("K132").Value = 5923.6362
$ws.Range("M132").Value = -3393.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 32255.445
$ws.Range("I7").Value = 39757
$ws.Range("K7").Value = 39757
$ws.Range("M7").Value = -39645

$ws.Range("H22").Value = 1391
$ws.Range("I22").Value = 1489.5
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 1489.5
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -1194.5
$ws.Range("N22").Value = -1390

$ws.Range("H27").Value = 1391
$ws.Range("I27").Value = 1489.5
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 1489.5
$ws.Range("L27").Value = 800
$ws.Range("M27").Value = -1382.5
$ws.Range("N27").Value = -1014

$ws.Range("H46").Value = 49846.332
$ws.Range("I46").Value = 62802.43
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 62802.43
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -62614.43
$ws.Range("N46").Value = -4876

$ws.Range("H126").Value = 32255.445
$ws.Range("I126").Value = 39757
$ws.Range("K126").Value = 119271
$ws.Range("M126").Value = -116801

$ws.Range("H136").Value = 5498.2856
$ws.Range("I136").Value = 3166.3333
$ws.Range("K136").Value = 9498.999899999999
$ws.Range("M136").Value = -6948.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 10000
$ws.Range("J28").Value = 10000
$ws.Range("L28").Value = 10000
$ws.Range("N28").Value = -10696

$ws.Range("H62").Value = 162509.5
$ws.Range("I62").Value = 8367
$ws.Range("K62").Value = 8367
$ws.Range("M62").Value = -7743

$ws.Range("H65").Value = 162509.5
$ws.Range("I65").Value = 8367
$ws.Range("K65").Value = 41835
$ws.Range("M65").Value = -38715
